$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "98.688.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.81%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.437.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.30%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "259.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.92%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "671.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.16%  "

$ws.Range("E7").Value = "  +10.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.468"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +18.30%  "

$ws.Range("E9").Value = "  +22.89%  "

$ws.Range("E10").Value = "  -0.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.434.21"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.34%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.220"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +11.43%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "43.19"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +13.85%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000276"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +12.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +12.48%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "98.322.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.69%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.068.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +34.27%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.427.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.97%  "

$ws.Range("E20").Value = "  +15.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "538.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +13.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.58%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +15.99%  "

$ws.Range("E24").Value = "  +8.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.437"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +51.90%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +15.82%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "102.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +16.84%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.88%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.605.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.87%  "

$ws.Range("E30").Value = "  +15.62%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +20.79%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.201"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.74%  "

$ws.Range("E33").Value = "  +0.06%  "

$ws.Range("E34").Value = "  +1.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "30.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.39%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.564"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +24.40%  "

$ws.Range("E37").Value = "  +15.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.95"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +11.83%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.161"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.65%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "535.10"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.66%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +15.32%  "

$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0441"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +37.32%  "

$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.867"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.41%  "

$ws.Range("B45").Value = "MantraDAO"
$ws.Range("C45").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.03%  "

$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.91%  "

$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +14.33%  "

$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.06%  "

$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +16.79%  "

$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +18.12%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +18.16%  "
